# Re-create the author's edit on Data/any_data/count.xlsx:
#   - the "t_turn_off" value on row 2 was changed from 7:00 (0.291666...)
#     to 6:30 (0.270833...)
#   - the cursor/selection left on the sheet moved from B8 to F7
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the 0:00 -> t_turn_off interval; update t_turn_off (column B)
# from 7:00 AM to 6:30 AM.
$ws.Range("B2").Value = 0.27083333333333331

# Leave the selection where the author left it before saving.
$ws.Range("F7").Select()
